# Applies the updated weight-estimation figures to the JPAD Weights workbook.
# These cells hold externally-computed (JPAD tool) static values, not Excel
# formulas, so each touched cell is written explicitly with its new value.

$wb = $excel.ActiveWorkbook

# --- GLOBAL RESULTS ---------------------------------------------------
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C6").Value  = 57381.743724072745
$ws.Range("C7").Value  = 57946.24372407276
$ws.Range("C8").Value  = 52151.61935166548
$ws.Range("C12").Value = 45987.1770447921
$ws.Range("C13").Value = 45987.1770447921
$ws.Range("C14").Value = 33117.17704479211
$ws.Range("C15").Value = 32388.089753792105
$ws.Range("C16").Value = 31537.879753792113
$ws.Range("C20").Value = 562722.6770916779
$ws.Range("C21").Value = 568258.5310166781
$ws.Range("C22").Value = 511432.6779150102
$ws.Range("C26").Value = 450980.14976631035
$ws.Range("C27").Value = 450980.14976631035
$ws.Range("C28").Value = 324768.5642663104
$ws.Range("C29").Value = 317618.6603840253
$ws.Range("C30").Value = 309280.94848752534

# --- FUSELAGE -----------------------------------------------------------
$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C8").Value  = 5986.0
$ws.Range("D8").Value  = 19.607568885247893
$ws.Range("C9").Value  = 2662.0
$ws.Range("D9").Value  = -46.80999860131475
$ws.Range("C12").Value = 4551.166666666666
$ws.Range("D12").Value = -9.06214824731418

# --- WING -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("WING")
$ws.Range("C7").Value  = 5487.0
$ws.Range("D7").Value  = 70.98251846312053
$ws.Range("C8").Value  = 4287.0
$ws.Range("D8").Value  = 33.588856688791275
$ws.Range("C9").Value  = 2790.0
$ws.Range("D9").Value  = -13.059736374684475
$ws.Range("C11").Value = 5818.0
$ws.Range("D11").Value = 81.29693683587301
$ws.Range("C12").Value = 5017.0
$ws.Range("D12").Value = 56.33666760150824
$ws.Range("C13").Value = 4086.7142857142853
$ws.Range("D13").Value = 27.347676473599634

# --- NACELLES ---------------------------------------------------------
$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("C3").Value  = 642.6666666666665
$ws.Range("D3").Value  = -83.6071149202462
$ws.Range("C9").Value  = 392.0
$ws.Range("D9").Value  = -40.00612182430363
$ws.Range("C11").Value = 323.0
$ws.Range("D11").Value = -50.56626874808692
$ws.Range("C12").Value = 321.33333333333326
$ws.Range("C16").Value = 392.0
$ws.Range("D16").Value = -40.00612182430363
$ws.Range("C18").Value = 323.0
$ws.Range("D18").Value = -50.56626874808692
$ws.Range("C19").Value = 321.33333333333326

# --- LANDING GEARS ------------------------------------------------------
$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C6").Value = 2290.0
$ws.Range("D6").Value = 2.5480274058483974
$ws.Range("C7").Value = 2616.0
$ws.Range("D7").Value = 17.146567551833805
$ws.Range("C8").Value = 2265.0
$ws.Range("D8").Value = 1.4285074560028908
$ws.Range("C9").Value = 2203.25
$ws.Range("D9").Value = -1.3367068201155303
